$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Cells.Item(1, 1).Value = "strain"
$ws.Cells.Item(1, 2).Value = "phage"
$ws.Cells.Item(1, 3).Value = "mean log10 PFU/mL"
$ws.Cells.Item(1, 4).Value = "sd log10 PFU/mL"
$ws.Cells.Item(1, 5).Value = "n"
$ws.Cells.Item(1, 6).Value = "n_not_NA"

# Data rows: strain, phage, mean, sd, n, n_not_NA
$data = @(
    @("b52",             "α15",   9.99287898279482, 0.266940829737246, 3, 3),
    @("b52",             "α17",   9.82686059560972, 0.111328624652246, 3, 3),
    @("b52",             "α15.2", 9.65167130101268, 0.19920152107145,  3, 3),
    @("b52∆LPS",         "α15",   0,                 0,                 3, 3),
    @("b52∆LPS",         "α17",   9.76110040956785, 0.151237604291255, 3, 3),
    @("b52∆LPS",         "α15.2", 9.61219242494688, 0.122470672390741, 3, 3),
    @("b52∆Tsx",         "α15",   9.76110040956785, 0.151237604291255, 3, 3),
    @("b52∆Tsx",         "α17",   0,                 0,                 3, 3),
    @("b52∆Tsx",         "α15.2", 6.38612083069842, 1.24472729672886,  3, 3),
    @("b52∆Tsx_and_LPS", "α15",   0,                 0,                 3, 3),
    @("b52∆Tsx_and_LPS", "α17",   0,                 0,                 3, 3),
    @("b52∆Tsx_and_LPS", "α15.2", 0,                 0,                 3, 3)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
